# #CRM-1168 Add brand column in Partner panel - Pending Spares - Download file
#
# Adds a new "Brand" column (O) to the Spare Requested Parts template:
#   O1 = "Brand"            (header, same style as the other header cells)
#   O2 = "{spare:brands}"   (placeholder token, same style as the neighbouring
#                             placeholder cell N2 / N1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new header + placeholder values.
$ws.Range("O1").Value = "Brand"
$ws.Range("O2").Value = "{spare:brands}"

# Copy the formatting from the existing header cell (A1) onto the new
# header cell (O1) so it keeps the same bold / centered / shaded style as
# the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Copy the formatting from the neighbouring placeholder cell (N2) onto the
# new placeholder cell (O2).
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
